$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.7799968719482422
$ws.Range("E2").Value = 1095.783988649098
$ws.Range("F2").Value = 0.06209605675574132
$ws.Range("G2").Value = 0.04444679389226319
$ws.Range("H2").Value = 0.03781111209987692
$ws.Range("I2").Value = 0.03228484000150536
$ws.Range("J2").Value = 0.02988204431794883
$ws.Range("K2").Value = 0.02715927593703929
$ws.Range("L2").Value = 0.02616427191177966
$ws.Range("M2").Value = 0.02458161152219122
$ws.Range("N2").Value = 0.02311344715861738
$ws.Range("O2").Value = 0.02295280704897289
$ws.Range("P2").Value = 0.02184812553522507
$ws.Range("Q2").Value = 0.02184812553522507
$ws.Range("R2").Value = 0.02184812553522507
$ws.Range("S2").Value = 0.02184812553522507
$ws.Range("T2").Value = 0.0216235132688433
$ws.Range("U2").Value = 0.0216235132688433
$ws.Range("V2").Value = 0.02144527603329724
$ws.Range("W2").Value = 0.02140713464867345
$ws.Range("X2").Value = 0.02137493487443794
$ws.Range("Y2").Value = 0.02136031166957306
$ws.Range("C3").Value = 0.699998140335083
$ws.Range("E3").Value = 1091.927398875259
$ws.Range("F3").Value = 0.0630556718683309
$ws.Range("G3").Value = 0.04592129906679353
$ws.Range("H3").Value = 0.03469268956972425
$ws.Range("I3").Value = 0.03229567441405932
$ws.Range("J3").Value = 0.02860913943704322
$ws.Range("K3").Value = 0.02674055548079961
$ws.Range("L3").Value = 0.02580615760969677
$ws.Range("M3").Value = 0.02446133727127513
$ws.Range("N3").Value = 0.02425656279078497
$ws.Range("O3").Value = 0.02337533115624447
$ws.Range("P3").Value = 0.02337533115624447
$ws.Range("Q3").Value = 0.0233407253966823
$ws.Range("R3").Value = 0.02200599556437013
$ws.Range("S3").Value = 0.02200599556437013
$ws.Range("T3").Value = 0.02193398224591377
$ws.Range("U3").Value = 0.02193398224591377
$ws.Range("V3").Value = 0.02164551619027143
$ws.Range("W3").Value = 0.0212851344809992
$ws.Range("X3").Value = 0.0212851344809992
$ws.Range("Y3").Value = 0.0212851344809992
$ws.Range("C4").Value = 0.7010025978088379
$ws.Range("E4").Value = 1167.282008417462
$ws.Range("F4").Value = 0.06065856373672367
$ws.Range("G4").Value = 0.04471597018877559
$ws.Range("H4").Value = 0.03777756807659313
$ws.Range("I4").Value = 0.03290402458334895
$ws.Range("J4").Value = 0.03086538355026833
$ws.Range("K4").Value = 0.02876788314372378
$ws.Range("L4").Value = 0.02677525690960798
$ws.Range("M4").Value = 0.02603924891308698
$ws.Range("N4").Value = 0.02504577590697768
$ws.Range("O4").Value = 0.02389409206935945
$ws.Range("P4").Value = 0.02389409206935945
$ws.Range("Q4").Value = 0.02389409206935945
$ws.Range("R4").Value = 0.0236822047155036
$ws.Range("S4").Value = 0.02355368333869444
$ws.Range("T4").Value = 0.02316355808123921
$ws.Range("U4").Value = 0.02314860349709712
$ws.Range("V4").Value = 0.02284697242839967
$ws.Range("W4").Value = 0.02284697242839967
$ws.Range("X4").Value = 0.02282340691226473
$ws.Range("Y4").Value = 0.02275403525180237
$ws.Range("C5").Value = 0.7739944458007812
$ws.Range("E5").Value = 1146.45400172902
$ws.Range("F5").Value = 0.06547067673107672
$ws.Range("G5").Value = 0.04934709161591286
$ws.Range("H5").Value = 0.04006777310888992
$ws.Range("I5").Value = 0.03687434548044861
$ws.Range("J5").Value = 0.03099272781391694
$ws.Range("K5").Value = 0.02896291592900091
$ws.Range("L5").Value = 0.02689687722120588
$ws.Range("M5").Value = 0.0261289944452383
$ws.Range("N5").Value = 0.02469076320935572
$ws.Range("O5").Value = 0.02469076320935572
$ws.Range("P5").Value = 0.02456202772746407
$ws.Range("Q5").Value = 0.02393282578671791
$ws.Range("R5").Value = 0.0235445021335369
$ws.Range("S5").Value = 0.02333978513036256
$ws.Range("T5").Value = 0.0229249887472257
$ws.Range("U5").Value = 0.0228269150425391
$ws.Range("V5").Value = 0.02277296506224588
$ws.Range("W5").Value = 0.02264102343314251
$ws.Range("X5").Value = 0.02242234111638487
$ws.Range("Y5").Value = 0.02234803122278791
$ws.Range("C6").Value = 0.754997730255127
$ws.Range("E6").Value = 1055.601810166045
$ws.Range("F6").Value = 0.06262306869582172
$ws.Range("G6").Value = 0.04491817338467746
$ws.Range("H6").Value = 0.03854453593921026
$ws.Range("I6").Value = 0.03479276141217992
$ws.Range("J6").Value = 0.03027172918842988
$ws.Range("K6").Value = 0.0274735465781626
$ws.Range("L6").Value = 0.0264044878525281
$ws.Range("M6").Value = 0.02552265932413967
$ws.Range("N6").Value = 0.02373297848475171
$ws.Range("O6").Value = 0.0226141400971253
$ws.Range("P6").Value = 0.0226141400971253
$ws.Range("Q6").Value = 0.02245607166364165
$ws.Range("R6").Value = 0.02201749529668413
$ws.Range("S6").Value = 0.02164458726976945
$ws.Range("T6").Value = 0.0208334589023506
$ws.Range("U6").Value = 0.0208334589023506
$ws.Range("V6").Value = 0.02068458618210488
$ws.Range("W6").Value = 0.02068458618210488
$ws.Range("X6").Value = 0.02057779566033666
$ws.Range("Y6").Value = 0.02057703333657007
$ws.Range("C7").Value = 0.7050409317016602
$ws.Range("E7").Value = 1107.713103692309
$ws.Range("F7").Value = 0.06197420643632718
$ws.Range("G7").Value = 0.04376011311779698
$ws.Range("H7").Value = 0.03999598845296207
$ws.Range("I7").Value = 0.031577154673315
$ws.Range("J7").Value = 0.02982723773808174
$ws.Range("K7").Value = 0.02781028102945483
$ws.Range("L7").Value = 0.02735011738089303
$ws.Range("M7").Value = 0.02577747209232614
$ws.Range("N7").Value = 0.02490054530726134
$ws.Range("O7").Value = 0.02393746215941415
$ws.Range("P7").Value = 0.02323579351555687
$ws.Range("Q7").Value = 0.02300725094401265
$ws.Range("R7").Value = 0.02280070719718902
$ws.Range("S7").Value = 0.02268847321221036
$ws.Range("T7").Value = 0.02268847321221036
$ws.Range("U7").Value = 0.02230882000453971
$ws.Range("V7").Value = 0.02205522715834128
$ws.Range("W7").Value = 0.02194641725348693
$ws.Range("X7").Value = 0.02173077753957404
$ws.Range("Y7").Value = 0.02159284802519121
$ws.Range("C8").Value = 0.7409617900848389
$ws.Range("E8").Value = 1144.653912324813
$ws.Range("F8").Value = 0.06055767099792104
$ws.Range("G8").Value = 0.04627093779875122
$ws.Range("H8").Value = 0.04014663360471366
$ws.Range("I8").Value = 0.03463111624609938
$ws.Range("J8").Value = 0.0319246799991795
$ws.Range("K8").Value = 0.0304195850695074
$ws.Range("L8").Value = 0.0291056603370582
$ws.Range("M8").Value = 0.02686866013514362
$ws.Range("N8").Value = 0.02543500233102915
$ws.Range("O8").Value = 0.0248454306520004
$ws.Range("P8").Value = 0.02451107969441342
$ws.Range("Q8").Value = 0.02387227938586449
$ws.Range("R8").Value = 0.02321299492791523
$ws.Range("S8").Value = 0.02311721282459537
$ws.Range("T8").Value = 0.0225396573475062
$ws.Range("U8").Value = 0.0225396573475062
$ws.Range("V8").Value = 0.0225396573475062
$ws.Range("W8").Value = 0.02242196606264447
$ws.Range("X8").Value = 0.02237036141632104
$ws.Range("Y8").Value = 0.0223129417607176
$ws.Range("C9").Value = 0.7840027809143066
$ws.Range("E9").Value = 1108.541484884439
$ws.Range("F9").Value = 0.06151579842035126
$ws.Range("G9").Value = 0.04552087913755756
$ws.Range("H9").Value = 0.03518479785350093
$ws.Range("I9").Value = 0.0317608377021031
$ws.Range("J9").Value = 0.02891781259409814
$ws.Range("K9").Value = 0.02746889438109408
$ws.Range("L9").Value = 0.02523148348914845
$ws.Range("M9").Value = 0.02523148348914845
$ws.Range("N9").Value = 0.02468181657887227
$ws.Range("O9").Value = 0.02403640635484982
$ws.Range("P9").Value = 0.02367451369416748
$ws.Range("Q9").Value = 0.02272281363205516
$ws.Range("R9").Value = 0.02272281363205516
$ws.Range("S9").Value = 0.02262501465404992
$ws.Range("T9").Value = 0.022427279189855
$ws.Range("U9").Value = 0.02216858258821053
$ws.Range("V9").Value = 0.02196600433015844
$ws.Range("W9").Value = 0.02180714500354845
$ws.Range("X9").Value = 0.02165182074066428
$ws.Range("Y9").Value = 0.02160899580671421
$ws.Range("C10").Value = 0.7139980792999268
$ws.Range("E10").Value = 1058.103226300167
$ws.Range("F10").Value = 0.06151636816357051
$ws.Range("G10").Value = 0.04407211891806202
$ws.Range("H10").Value = 0.03804203689454911
$ws.Range("I10").Value = 0.03344441370197019
$ws.Range("J10").Value = 0.02956204483274667
$ws.Range("K10").Value = 0.0290459607246605
$ws.Range("L10").Value = 0.02647577374843785
$ws.Range("M10").Value = 0.02489441558537952
$ws.Range("N10").Value = 0.02489441558537952
$ws.Range("O10").Value = 0.02334800208243291
$ws.Range("P10").Value = 0.02259697272682696
$ws.Range("Q10").Value = 0.02221402025732742
$ws.Range("R10").Value = 0.02185108513856289
$ws.Range("S10").Value = 0.02133289316651324
$ws.Range("T10").Value = 0.02133289316651324
$ws.Range("U10").Value = 0.02099481858006461
$ws.Range("V10").Value = 0.02099481858006461
$ws.Range("W10").Value = 0.02092331083508847
$ws.Range("X10").Value = 0.0206918911217712
$ws.Range("Y10").Value = 0.02062579388499351
$ws.Range("C11").Value = 0.6939990520477295
$ws.Range("E11").Value = 1068.888232259978
$ws.Range("F11").Value = 0.0628823000715748
$ws.Range("G11").Value = 0.05145934003302075
$ws.Range("H11").Value = 0.03867223923643547
$ws.Range("I11").Value = 0.03447745482925521
$ws.Range("J11").Value = 0.02955501028024562
$ws.Range("K11").Value = 0.02712016958365799
$ws.Range("L11").Value = 0.0247689457158031
$ws.Range("M11").Value = 0.0247689457158031
$ws.Range("N11").Value = 0.02407480745971677
$ws.Range("O11").Value = 0.02288523610028501
$ws.Range("P11").Value = 0.02284272842364613
$ws.Range("Q11").Value = 0.02244444093605468
$ws.Range("R11").Value = 0.02190627752422604
$ws.Range("S11").Value = 0.02190627752422604
$ws.Range("T11").Value = 0.02166007130199268
$ws.Range("U11").Value = 0.02127498524119347
$ws.Range("V11").Value = 0.02127498524119347
$ws.Range("W11").Value = 0.02106115400611442
$ws.Range("X11").Value = 0.02103624107515012
$ws.Range("Y11").Value = 0.02083602791929782

Write-Output "Applied 220 cell updates"